# Auto-generated edit script applying Sheets update per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 403.5
$ws.Range("I38").Value = 45.6
$ws.Range("K38").Value = 136.8
$ws.Range("M38").Value = 235.2

$ws.Range("H88").Value = 2388.4167
$ws.Range("I88").Value = 1968.2
$ws.Range("J88").Value = 2688.5715
$ws.Range("K88").Value = 1968.2
$ws.Range("L88").Value = 2688.5715
$ws.Range("M88").Value = -1562.2
$ws.Range("N88").Value = -3500.5715

$ws.Range("H91").Value = 2388.4167
$ws.Range("I91").Value = 1968.2
$ws.Range("J91").Value = 2688.5715
$ws.Range("K91").Value = 1968.2
$ws.Range("L91").Value = 2688.5715
$ws.Range("M91").Value = -564.2
$ws.Range("N91").Value = -5496.5715

$ws.Range("H98").Value = 611.5238000000001
$ws.Range("I98").Value = 627.3333
$ws.Range("J98").Value = 516.6667
$ws.Range("K98").Value = 627.3333
$ws.Range("L98").Value = 516.6667
$ws.Range("M98").Value = 870.6667
$ws.Range("N98").Value = -3512.6667

$ws.Range("H122").Value = 611.5238000000001
$ws.Range("I122").Value = 627.3333
$ws.Range("J122").Value = 516.6667
$ws.Range("K122").Value = 1881.9999
$ws.Range("L122").Value = 1550.0001
$ws.Range("M122").Value = 568.0001
$ws.Range("N122").Value = -6450.0001

$ws.Range("H129").Value = 357878.16
$ws.Range("I129").Value = 274.25
$ws.Range("J129").Value = 417478.78
$ws.Range("K129").Value = 822.75
$ws.Range("L129").Value = 1252436.34
$ws.Range("M129").Value = 4177.25
$ws.Range("N129").Value = -1262436.34

$ws.Range("H132").Value = 3752.4348
$ws.Range("I132").Value = 4426.6313
$ws.Range("J132").Value = 550
$ws.Range("K132").Value = 13279.8939
$ws.Range("L132").Value = 1650
$ws.Range("M132").Value = -10749.8939
$ws.Range("N132").Value = -6710

$ws.Range("H135").Value = 31260390
$ws.Range("I135").Value = 1022.4545
$ws.Range("K135").Value = 9202.0905
$ws.Range("M135").Value = -6667.0905

$ws.Range("H137").Value = 2151.4614
$ws.Range("I137").Value = 1886.9
$ws.Range("J137").Value = 3033.3333
$ws.Range("K137").Value = 5660.700000000001
$ws.Range("L137").Value = 9099.999899999999
$ws.Range("M137").Value = -3110.700000000001
$ws.Range("N137").Value = -14199.9999

$ws.Range("H138").Value = 2260.013
$ws.Range("I138").Value = 2070.0625
$ws.Range("J138").Value = 2309.0322
$ws.Range("K138").Value = 6210.1875
$ws.Range("L138").Value = 6927.096600000001
$ws.Range("M138").Value = -1070.1875
$ws.Range("N138").Value = -17207.0966

$ws.Range("H141").Value = 2152.647
$ws.Range("I141").Value = 1529.091
$ws.Range("K141").Value = 4587.272999999999
$ws.Range("M141").Value = 592.7270000000008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7432.2925
$ws.Range("I32").Value = 5390.359
$ws.Range("K32").Value = 5390.359
$ws.Range("M32").Value = -5103.359

$ws.Range("H132").Value = 19328.035
$ws.Range("I132").Value = 2000.7142
$ws.Range("J132").Value = 64812.25
$ws.Range("K132").Value = 6002.142599999999
$ws.Range("L132").Value = 194436.75
$ws.Range("M132").Value = -3472.142599999999
$ws.Range("N132").Value = -199496.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1443.7
$ws.Range("I99").Value = 1128
$ws.Range("J99").Value = 2030
$ws.Range("K99").Value = 1128
$ws.Range("L99").Value = 2030
$ws.Range("M99").Value = 370
$ws.Range("N99").Value = -5026

$ws.Range("H134").Value = 3638.9707
$ws.Range("I134").Value = 3890.8333
$ws.Range("J134").Value = 1750
$ws.Range("K134").Value = 11672.4999
$ws.Range("L134").Value = 5250
$ws.Range("M134").Value = -9137.499899999999
$ws.Range("N134").Value = -10320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 18084.848
$ws.Range("I132").Value = 24953.682
$ws.Range("K132").Value = 74861.046
$ws.Range("M132").Value = -72331.046

$ws.Range("H134").Value = 1234.2
$ws.Range("I134").Value = 1192.2307
$ws.Range("J134").Value = 1507
$ws.Range("K134").Value = 3576.6921
$ws.Range("L134").Value = 4521
$ws.Range("M134").Value = -1041.6921
$ws.Range("N134").Value = -9591

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 766.51
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 766.51
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2299.53
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12379.53

$ws.Range("H141").Value = 3485.1333
$ws.Range("I141").Value = 3106.5833
$ws.Range("J141").Value = 4999.3335
$ws.Range("K141").Value = 9319.749899999999
$ws.Range("L141").Value = 14998.0005
$ws.Range("M141").Value = -4139.749899999999
$ws.Range("N141").Value = -25358.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7000
$ws.Range("J43").Value = 7000
$ws.Range("L43").Value = 7000
$ws.Range("N43").Value = -7302

$ws.Range("H57").Value = 27278
$ws.Range("I57").Value = 26200
$ws.Range("J57").Value = 27996.666
$ws.Range("K57").Value = 26200
$ws.Range("L57").Value = 27996.666
$ws.Range("M57").Value = -25380
$ws.Range("N57").Value = -29636.666

$ws.Range("H97").Value = 1097.909
$ws.Range("I97").Value = 1008.55554
$ws.Range("K97").Value = 1008.55554
$ws.Range("M97").Value = -512.55554

$ws.Range("H102").Value = 1241.1177
$ws.Range("I102").Value = 1385.9166
$ws.Range("J102").Value = 893.6
$ws.Range("K102").Value = 1385.9166
$ws.Range("L102").Value = 893.6
$ws.Range("M102").Value = 236.0834
$ws.Range("N102").Value = -4137.6

$ws.Range("H132").Value = 30037.8
$ws.Range("I132").Value = 6058.3076
$ws.Range("J132").Value = 74571.14
$ws.Range("K132").Value = 18174.9228
$ws.Range("L132").Value = 223713.42
$ws.Range("M132").Value = -15644.9228
$ws.Range("N132").Value = -228773.42

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5456.353
$ws.Range("I40").Value = 3761.7778
$ws.Range("J40").Value = 7362.75
$ws.Range("K40").Value = 3761.7778
$ws.Range("L40").Value = 7362.75
$ws.Range("M40").Value = -3625.7778
$ws.Range("N40").Value = -7634.75

$ws.Range("H55").Value = 78.72727
$ws.Range("I55").Value = 80
$ws.Range("J55").Value = 78.25
$ws.Range("K55").Value = 80
$ws.Range("L55").Value = 78.25
$ws.Range("M55").Value = 93
$ws.Range("N55").Value = -424.25

$ws.Range("H100").Value = 2888.3684
$ws.Range("I100").Value = 1396.3334
$ws.Range("J100").Value = 3168.125
$ws.Range("K100").Value = 1396.3334
$ws.Range("L100").Value = 3168.125
$ws.Range("M100").Value = -855.3334
$ws.Range("N100").Value = -4250.125

$ws.Range("H132").Value = 1996.6552
$ws.Range("I132").Value = 1518.4546
$ws.Range("J132").Value = 2288.889
$ws.Range("K132").Value = 4555.3638
$ws.Range("L132").Value = 6866.667
$ws.Range("M132").Value = -2025.3638
$ws.Range("N132").Value = -11926.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 1711.25
$ws.Range("J6").Value = 1946.6666
$ws.Range("L6").Value = 1946.6666
$ws.Range("N6").Value = -2176.6666

$ws.Range("H46").Value = 27126.334
$ws.Range("J46").Value = 27126.334
$ws.Range("L46").Value = 27126.334
$ws.Range("N46").Value = -27588.334

$ws.Range("H132").Value = 922.9535
$ws.Range("I132").Value = 648.05884
$ws.Range("K132").Value = 1944.17652
$ws.Range("M132").Value = 585.82348

$ws.Range("H134").Value = 27126.334
$ws.Range("J134").Value = 27126.334
$ws.Range("L134").Value = 81379.00199999999
$ws.Range("N134").Value = -86449.00199999999
